$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("Sheet1")

# Row 5: mark risk as Resolved ("R"), record the mitigation date, and
# update the "How" note with the new text.
$ws.Range("F5").Value = 45043
$ws.Range("G5").Value = "R"
$ws.Range("H5").Value = "I finished early so I can now earn 5 bonus points with early demo"

# Update the active selection to match the saved view state.
$ws.Range("H13").Select()
